# Applies the weekly Espinaca data refresh:
# - a new record is inserted at row 283 (pushing rows 283-299 down to 284-300)
# - row 300 is a brand-new row carrying what used to be row 299's data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 283 becomes a new record (scaffold columns A,B,C,E,F,G,H,I,N,O,Q,R stay the same
# as every other Espinaca row in this block, so only the data columns change).
$ws.Range("D283").Value = 44783
$ws.Range("J283").Value = 2600
$ws.Range("K283").Value = 500
$ws.Range("L283").Value = 600
$ws.Range("M283").Value = 550
$ws.Range("P283").Value = 1100

# Rows 284-299 each take on the data that used to live one row above them
# row 284 <- old row 283
$ws.Range("D284").Value = 44377
$ws.Range("J284").Value = 3120
$ws.Range("K284").Value = 400
$ws.Range("L284").Value = 500
$ws.Range("M284").Value = 450
$ws.Range("P284").Value = 900

# row 285 <- old row 284
$ws.Range("D285").Value = 44512
$ws.Range("J285").Value = 3340
$ws.Range("K285").Value = 400
$ws.Range("L285").Value = 500
$ws.Range("M285").Value = 450
$ws.Range("P285").Value = 900

# row 286 <- old row 285
$ws.Range("D286").Value = 44497
$ws.Range("J286").Value = 3000
$ws.Range("K286").Value = 400
$ws.Range("L286").Value = 500
$ws.Range("M286").Value = 450
$ws.Range("P286").Value = 900

# row 287 <- old row 286
$ws.Range("D287").Value = 44557
$ws.Range("J287").Value = 2400
$ws.Range("K287").Value = 400
$ws.Range("L287").Value = 500
$ws.Range("M287").Value = 450
$ws.Range("P287").Value = 900

# row 288 <- old row 287
$ws.Range("D288").Value = 44357
$ws.Range("J288").Value = 3000
$ws.Range("K288").Value = 450
$ws.Range("L288").Value = 500
$ws.Range("M288").Value = 475
$ws.Range("P288").Value = 950

# row 289 <- old row 288
$ws.Range("D289").Value = 44279
$ws.Range("J289").Value = 3200
$ws.Range("K289").Value = 400
$ws.Range("L289").Value = 500
$ws.Range("M289").Value = 450
$ws.Range("P289").Value = 900

# row 290 <- old row 289
$ws.Range("D290").Value = 44517
$ws.Range("J290").Value = 3000
$ws.Range("K290").Value = 450
$ws.Range("L290").Value = 500
$ws.Range("M290").Value = 475
$ws.Range("P290").Value = 950

# row 291 <- old row 290
$ws.Range("D291").Value = 44547
$ws.Range("J291").Value = 3320
$ws.Range("K291").Value = 400
$ws.Range("L291").Value = 500
$ws.Range("M291").Value = 450
$ws.Range("P291").Value = 900

# row 292 <- old row 291
$ws.Range("D292").Value = 44321
$ws.Range("J292").Value = 3060
$ws.Range("K292").Value = 450
$ws.Range("L292").Value = 500
$ws.Range("M292").Value = 475
$ws.Range("P292").Value = 950

# row 293 <- old row 292
$ws.Range("D293").Value = 44438
$ws.Range("J293").Value = 3460
$ws.Range("K293").Value = 400
$ws.Range("L293").Value = 500
$ws.Range("M293").Value = 450
$ws.Range("P293").Value = 900

# row 294 <- old row 293
$ws.Range("D294").Value = 44657
$ws.Range("J294").Value = 2460
$ws.Range("K294").Value = 500
$ws.Range("L294").Value = 600
$ws.Range("M294").Value = 550
$ws.Range("P294").Value = 1100

# row 295 <- old row 294
$ws.Range("D295").Value = 44391
$ws.Range("J295").Value = 3100
$ws.Range("K295").Value = 450
$ws.Range("L295").Value = 500
$ws.Range("M295").Value = 475
$ws.Range("P295").Value = 950

# row 296 <- old row 295
$ws.Range("D296").Value = 44186
$ws.Range("J296").Value = 2400
$ws.Range("K296").Value = 400
$ws.Range("L296").Value = 500
$ws.Range("M296").Value = 450
$ws.Range("P296").Value = 900

# row 297 <- old row 296
$ws.Range("D297").Value = 44189
$ws.Range("J297").Value = 2800
$ws.Range("K297").Value = 400
$ws.Range("L297").Value = 500
$ws.Range("M297").Value = 450
$ws.Range("P297").Value = 900

# row 298 <- old row 297
$ws.Range("D298").Value = 44609
$ws.Range("J298").Value = 2000
$ws.Range("K298").Value = 500
$ws.Range("L298").Value = 600
$ws.Range("M298").Value = 550
$ws.Range("P298").Value = 1100

# row 299 <- old row 298
$ws.Range("D299").Value = 44489
$ws.Range("J299").Value = 3000
$ws.Range("K299").Value = 450
$ws.Range("L299").Value = 500
$ws.Range("M299").Value = 475
$ws.Range("P299").Value = 950

# Row 300 is brand new: it carries what used to be row 299's full record
$ws.Range("A300").Value = 8
$ws.Range("B300").Value = "Terminal La Palmera de La Serena"
$ws.Range("C300").Value = "Coquimbo"
$ws.Range("D300").Value = 44358
$ws.Range("E300").Value = 4
$ws.Range("F300").Value = 100112012
$ws.Range("G300").Value = "Espinaca"
$ws.Range("H300").Value = "Sin especificar"
$ws.Range("I300").Value = "Primera"
$ws.Range("J300").Value = 3360
$ws.Range("K300").Value = 400
$ws.Range("L300").Value = 500
$ws.Range("M300").Value = 450
$ws.Range("N300").Value = "`$/atado 300 a 500 gramos"
$ws.Range("O300").Value = "Provincia del Elquí"
$ws.Range("P300").Value = 900
$ws.Range("Q300").Value = 0.5
$ws.Range("R300").Value = "Hortaliza"
$ws.Range("D300").NumberFormat = "YYYY-MM-DD HH:MM:SS"
